# Reorders the "Recorded By" list in column G: for any cell whose value is a
# comma-separated list of names/emails that ends with "System", the first
# entry is moved to the end of the list (left-rotation by one element).
# Cells that don't end with "System" (e.g. single names, or lists without
# "System") are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($null -eq $val) { continue }

    $text = [string]$val
    if ($text -notmatch ",") { continue }

    $parts = $text -split ",\s*"
    if ($parts.Count -lt 2) { continue }

    $last = $parts[$parts.Count - 1].Trim()
    if ($last -ne "System") { continue }

    $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
    $cell.Value2 = [string]::Join(", ", $rotated)
}
